$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B29 currently stores "4" as an inline/text string; convert it to a true number
$ws.Range("B29").Value = 4

# Add new row 30 data
$ws.Range("A30").Value = "Sunsi Wu"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "5"
$ws.Range("C30").Value = "insightful"
$ws.Range("D30").Value = "APC"
$ws.Range("E30").Value = "OTH"
$ws.Range("F30").Value = "7f314748-ac5a-4a11-8786-6125314f9d6d"
$ws.Range("G30").Value = "Sy2ogebAW_annotated.xlsx"
$ws.Range("H30").Value = "We would like to thank all reviewers for their detailed and insightful feedback."
